# chore: adapt column header formatting to respective input file names
# Rename the "_old"/"_new" column header suffixes to the respective
# format-version suffixes ("_FV2210" / "_FV2304"), turn the sheet's data
# range into a real Excel Table (ListObject) and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row (row 1) ------------------------------------
# Columns A:J used the "_old" suffix, column K is the literal "diff"
# column (left untouched) and columns L:U used the "_new" suffix.
$baseNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $leftCol = $i + 1          # columns 1..10  -> A..J
    $rightCol = $i + 12        # columns 12..21 -> L..U
    $ws.Cells.Item(1, $leftCol).Value = "$($baseNames[$i])_FV2210"
    $ws.Cells.Item(1, $rightCol).Value = "$($baseNames[$i])_FV2304"
}

# --- 2. Turn the data range into an Excel Table -----------------------
$tableRange = $ws.Range("A1:U56")
$lo = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$lo.Name = "Table1"

# --- 3. Freeze the header row ------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Host "Header renaming, table creation and freeze panes applied."
